# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.734.32"
$ws.Range("E2").Value = "  +3.09%  "

# Row 3
$ws.Range("D3").Value = "1.692.99"
$ws.Range("E3").Value = "  +3.60%  "

# Row 4
$ws.Range("D4").Value = "'0.9984"
$ws.Range("E4").Value = "  -0.46%  "

# Row 5
$ws.Range("D5").Value = "'218.34"
$ws.Range("E5").Value = "  +4.33%  "

# Row 6
$ws.Range("D6").Value = "'0.5374"
$ws.Range("E6").Value = "  +4.04%  "

# Row 7
$ws.Range("D7").Value = "'0.9984"
$ws.Range("E7").Value = "  -0.49%  "

# Row 8
$ws.Range("D8").Value = "'0.2696"
$ws.Range("E8").Value = "  +5.67%  "

# Row 9
$ws.Range("D9").Value = "'0.06460"
$ws.Range("E9").Value = "  +3.86%  "

# Row 10
$ws.Range("D10").Value = "'21.66"
$ws.Range("E10").Value = "  +7.23%  "

# Row 11
$ws.Range("D11").Value = "'0.07810"
$ws.Range("E11").Value = "  +3.53%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.693.00"
$ws.Range("E12").Value = "  +3.44%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.530"
$ws.Range("E13").Value = "  +4.29%  "

# Row 14
$ws.Range("D14").Value = "'0.5674"
$ws.Range("E14").Value = "  +5.04%  "

# Row 15
$ws.Range("D15").Value = "0.0₅8526"
$ws.Range("E15").Value = "  +7.99%  "

# Row 16
$ws.Range("D16").Value = "'66.56"
$ws.Range("E16").Value = "  +3.50%  "

# Row 17
$ws.Range("D17").Value = "26.726.85"
$ws.Range("E17").Value = "  +2.95%  "

# Row 18
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  -0.14%  "

# Row 19
$ws.Range("D19").Value = "'4.835"
$ws.Range("E19").Value = "  +4.79%  "

# Row 20
$ws.Range("D20").Value = "'196.86"
$ws.Range("E20").Value = "  +6.96%  "

# Row 21
$ws.Range("D21").Value = "'10.46"
$ws.Range("E21").Value = "  +5.08%  "

# Row 22
$ws.Range("D22").Value = "'6.422"
$ws.Range("E22").Value = "  +5.97%  "

# Row 23
$ws.Range("D23").Value = "'0.9999"
$ws.Range("E23").Value = "  -0.39%  "

# Row 24
$ws.Range("D24").Value = "'143.24"
$ws.Range("E24").Value = "  -1.32%  "

# Row 25
$ws.Range("D25").Value = "'0.1285"
$ws.Range("E25").Value = "  +8.41%  "

# Row 26
$ws.Range("D26").Value = "'7.508"
$ws.Range("E26").Value = "  +2.91%  "

# Row 27
$ws.Range("D27").Value = "'16.30"
$ws.Range("E27").Value = "  +5.65%  "

# Row 28
$ws.Range("D28").Value = "'1.420"

# Row 29
$ws.Range("D29").Value = "'0.06196"
$ws.Range("E29").Value = "  +4.90%  "

# Row 30
$ws.Range("D30").Value = "'1.282"
$ws.Range("E30").Value = "  +3.32%  "

# Row 31
$ws.Range("D31").Value = "'3.622"
$ws.Range("E31").Value = "  +9.01%  "

# Row 32
$ws.Range("D32").Value = "'3.482"
$ws.Range("E32").Value = "  +4.55%  "

# Row 33
$ws.Range("D33").Value = "'1.718"
$ws.Range("E33").Value = "  +7.52%  "

# Row 34
$ws.Range("D34").Value = "'1.022"
$ws.Range("E34").Value = "  +6.00%  "

# Row 35
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").Value = "'2.797"
$ws.Range("E35").Value = "  +2.25%  "

# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.416"
$ws.Range("E36").Value = "  +1.47%  "

# Row 37
$ws.Range("D37").Value = "'0.5750"
$ws.Range("E37").Value = "  -0.98%  "

# Row 38
$ws.Range("D38").Value = "'0.01652"
$ws.Range("E38").Value = "  +3.83%  "

# Row 39
$ws.Range("D39").Value = "'5.975"
$ws.Range("E39").Value = "  +5.43%  "

# Row 40
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.075.18"
$ws.Range("E40").Value = "  +4.03%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8676"
$ws.Range("E41").Value = "  +3.49%  "

# Row 42
$ws.Range("D42").Value = "'0.9995"
$ws.Range("E42").Value = "  -0.26%  "

# Row 43
$ws.Range("D43").Value = "'100.48"
$ws.Range("E43").Value = "  +0.99%  "

# Row 44
$ws.Range("D44").Value = "1.835.96"
$ws.Range("E44").Value = "  +2.66%  "

# Row 45
$ws.Range("D45").Value = "'57.59"
$ws.Range("E45").Value = "  +6.24%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈106"
$ws.Range("E46").Value = "  -1.71%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'8.173"
$ws.Range("E47").Value = "  +3.08%  "

# Row 48
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "  +0.30%  "

# Row 49
$ws.Range("D49").Value = "'6.121"
$ws.Range("E49").Value = "  +6.27%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05218"
$ws.Range("E50").Value = "  +0.64%  "

# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.4255"
$ws.Range("E51").Value = "  +0.58%  "
